# Applies odds/score updates scraped for Jogos_da_Semana_FlashScore_2025-03-07.xlsx
# Updates specific cells in rows 8, 19, 25, 51, 52 and 66 on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("G8").Value2 = 2.55
$ws.Range("I8").Value2 = 2.7
$ws.Range("J8").Value2 = 3.4
$ws.Range("K8").Value2 = 2
$ws.Range("L8").Value2 = 3.5
$ws.Range("M8").Value2 = 1.07
$ws.Range("N8").Value2 = 8.5
$ws.Range("O8").Value2 = 1.4
$ws.Range("P8").Value2 = 2.75
$ws.Range("Q8").Value2 = 2.2
$ws.Range("R8").Value2 = 1.62
$ws.Range("S8").Value2 = 4
$ws.Range("T8").Value2 = 1.22
$ws.Range("U8").Value2 = 1.5
$ws.Range("V8").Value2 = 2.5
$ws.Range("W8").Value2 = 1.95
$ws.Range("X8").Value2 = 1.8
$ws.Range("Y8").Value2 = 7
$ws.Range("AC8").Value2 = 23
$ws.Range("AD8").Value2 = 34
$ws.Range("AE8").Value2 = 8.5
$ws.Range("AG8").Value2 = 17
$ws.Range("AI8").Value2 = 401
$ws.Range("AJ8").Value2 = 7.5
$ws.Range("AK8").Value2 = 12
$ws.Range("AN8").Value2 = 23
$ws.Range("AO8").Value2 = 34
$ws.Range("AR8").Value2 = 3.25
$ws.Range("AS8").Value2 = 1.34
# Row 19
$ws.Range("S19").Value2 = 6.5
$ws.Range("T19").Value2 = 1.11
# Row 25
$ws.Range("G25").Value2 = 1.85
$ws.Range("H25").Value2 = 3.4
$ws.Range("I25").Value2 = 4.2
$ws.Range("Q25").Value2 = 2.15
$ws.Range("R25").Value2 = 1.67
$ws.Range("Z25").Value2 = 8
$ws.Range("AN25").Value2 = 41
# Row 51
$ws.Range("G51").Value2 = 2.35
$ws.Range("I51").Value2 = 3
$ws.Range("L51").Value2 = 3.5
$ws.Range("Q51").Value2 = 1.85
$ws.Range("R51").Value2 = 2
$ws.Range("AA51").Value2 = 9.5
$ws.Range("AD51").Value2 = 26
$ws.Range("AE51").Value2 = 11
$ws.Range("AG51").Value2 = 13
$ws.Range("AK51").Value2 = 15
$ws.Range("AM51").Value2 = 29
# Row 52
$ws.Range("G52").Value2 = 1.3
$ws.Range("H52").Value2 = 5.5
$ws.Range("I52").Value2 = 8.5
$ws.Range("J52").Value2 = 1.8
$ws.Range("K52").Value2 = 2.5
$ws.Range("L52").Value2 = 8
$ws.Range("O52").Value2 = 1.2
$ws.Range("P52").Value2 = 4.33
$ws.Range("Q52").Value2 = 1.65
$ws.Range("R52").Value2 = 2.2
$ws.Range("S52").Value2 = 2.63
$ws.Range("T52").Value2 = 1.44
$ws.Range("U52").Value2 = 1.3
$ws.Range("V52").Value2 = 3.4
$ws.Range("W52").Value2 = 2.05
$ws.Range("X52").Value2 = 1.7
$ws.Range("Z52").Value2 = 6.5
$ws.Range("AB52").Value2 = 8
$ws.Range("AC52").Value2 = 11
$ws.Range("AE52").Value2 = 15
$ws.Range("AF52").Value2 = 10
$ws.Range("AG52").Value2 = 23
$ws.Range("AH52").Value2 = 67
$ws.Range("AI52").Value2 = 401
$ws.Range("AJ52").Value2 = 21
$ws.Range("AL52").Value2 = 26
$ws.Range("AM52").Value2 = 101
$ws.Range("AN52").Value2 = 67
# Row 66
$ws.Range("AR66").Value2 = 1.95
$ws.Range("AS66").Value2 = 1.9
